$d = $word.ActiveDocument

$replacements = @(
    @("2024-04-16 Tuesday", "2024-04-17 Wednesday"),
    @("857÷3=", "266÷7="),
    @("373÷2=", "835÷4="),
    @("955÷3=", "191÷9="),
    @("417÷9=", "854÷4="),
    @("213÷9=", "426÷8="),
    @("805÷8=", "731÷4="),
    @("846÷7=", "157÷5="),
    @("952÷3=", "998÷7="),
    @("921÷4=", "893÷6="),
    @("501÷8=", "365÷9="),
    @("894÷3=", "469÷3="),
    @("567÷3=", "560÷4="),
    @("458÷8=", "621÷4="),
    @("108÷8=", "112÷8="),
    @("117÷2=", "719÷5="),
    @("582÷9=", "224÷9="),
    @("540÷4=", "466÷6="),
    @("422÷9=", "838÷7="),
    @("239÷3=", "154÷6="),
    @("652÷7=", "421÷7="),
    @("344÷2=", "659÷9="),
    @("420÷7=", "124÷5="),
    @("600÷4=", "153÷2="),
    @("702÷8=", "119÷4="),
    @("877÷4=", "142÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
